# Auto-generated Excel COM-interop edit script
# Applies cell-value updates to the Odin_Profits-derived sheets (ALC, ARM, BSM, CRP, GSM, LTW)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 3623.9473  # was 2688.423
$ws.Range("I5").Value = 990.46155  # was 734.44446
$ws.Range("J5").Value = 9329.833000000001  # was 7084.875
$ws.Range("K5").Value = 990.46155  # was 734.44446
$ws.Range("L5").Value = 9329.833000000001  # was 7084.875
$ws.Range("M5").Value = -875.46155  # was -619.44446
$ws.Range("N5").Value = -9559.833000000001  # was -7314.875

$ws.Range("H33").Value = 422.6  # was 442.47058
$ws.Range("I33").Value = 417.33334  # was 442.2143
$ws.Range("K33").Value = 417.33334  # was 442.2143
$ws.Range("M33").Value = -188.33334  # was -213.2143

$ws.Range("H80").Value = 456  # was 516.5454999999999
$ws.Range("I80").Value = 461  # was 544
$ws.Range("J80").Value = 444.75  # was 468.5
$ws.Range("K80").Value = 1383  # was 1632
$ws.Range("L80").Value = 1334.25  # was 1405.5
$ws.Range("M80").Value = -385  # was -634
$ws.Range("N80").Value = -3330.25  # was -3401.5

$ws.Range("H83").Value = 456  # was 516.5454999999999
$ws.Range("I83").Value = 461  # was 544
$ws.Range("J83").Value = 444.75  # was 468.5
$ws.Range("K83").Value = 4149  # was 4896
$ws.Range("L83").Value = 4002.75  # was 4216.5
$ws.Range("M83").Value = 843  # was 96
$ws.Range("N83").Value = -13986.75  # was -14200.5

$ws.Range("H92").Value = 747.2692  # was 749.1923
$ws.Range("J92").Value = 1524.25  # was 1536.75
$ws.Range("L92").Value = 1524.25  # was 1536.75
$ws.Range("N92").Value = -4020.25  # was -4032.75

$ws.Range("H103").Value = 1999.5  # was 0
$ws.Range("J103").Value = 1999.5  # was 0
$ws.Range("L103").Value = 5998.5  # was 0
$ws.Range("N103").Value = -7170.5  # was None

$ws.Range("H107").Value = 783.8823  # was 1036.5
$ws.Range("I107").Value = 921.2143  # was 1214.8
$ws.Range("J107").Value = 143  # was 145
$ws.Range("K107").Value = 921.2143  # was 1214.8
$ws.Range("L107").Value = 143  # was 145
$ws.Range("M107").Value = 998.7857  # was 705.2
$ws.Range("N107").Value = -3983  # was -3985

$ws.Range("H113").Value = 3727.389  # was 3829.1177
$ws.Range("I113").Value = 2998.8333  # was 3089.818
$ws.Range("K113").Value = 2998.8333  # was 3089.818
$ws.Range("M113").Value = 255.1667000000002  # was 164.1819999999998

$ws.Range("H129").Value = 2671.7693  # was 2678.2
$ws.Range("J129").Value = 3433.5715  # was 4021
$ws.Range("L129").Value = 10300.7145  # was 12063
$ws.Range("N129").Value = -20300.7145  # was -22063

$ws.Range("H137").Value = 18085.3  # was 14330.154
$ws.Range("I137").Value = 50000  # was 25850
$ws.Range("J137").Value = 14539.223  # was 12235.637
$ws.Range("K137").Value = 150000  # was 77550
$ws.Range("L137").Value = 43617.669  # was 36706.911
$ws.Range("M137").Value = -147450  # was -75000
$ws.Range("N137").Value = -48717.669  # was -41806.911

$ws.Range("H138").Value = 3544.7737  # was 3548.17
$ws.Range("I138").Value = 1048  # was 1066
$ws.Range("K138").Value = 3144  # was 3198
$ws.Range("M138").Value = 1996  # was 1942

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 350  # was 600.5
$ws.Range("I4").Value = 200  # was 600.5
$ws.Range("J4").Value = 500  # was 0
$ws.Range("K4").Value = 200  # was 600.5
$ws.Range("L4").Value = 500  # was 0
$ws.Range("M4").Value = -84  # was -484.5
$ws.Range("N4").Value = -732  # was None

$ws.Range("H5").Value = 820  # was 708.4
$ws.Range("I5").Value = 244.5  # was 260.5
$ws.Range("J5").Value = 1395.5  # was 2500
$ws.Range("K5").Value = 244.5  # was 260.5
$ws.Range("L5").Value = 1395.5  # was 2500
$ws.Range("M5").Value = -132.5  # was -148.5
$ws.Range("N5").Value = -1619.5  # was -2724

$ws.Range("H32").Value = 789.431  # was 819.1070999999999
$ws.Range("I32").Value = 503.63635  # was 524.2075
$ws.Range("K32").Value = 503.63635  # was 524.2075
$ws.Range("M32").Value = -216.63635  # was -237.2075

$ws.Range("H61").Value = 4846.7036  # was 4597.8623
$ws.Range("I61").Value = 5088.533  # was 4635.5884
$ws.Range("K61").Value = 5088.533  # was 4635.5884
$ws.Range("M61").Value = -4876.533  # was -4423.5884

$ws.Range("H88").Value = 6378.1113  # was 6429.778
$ws.Range("I88").Value = 963.3333  # was 1144.2
$ws.Range("J88").Value = 9085.5  # was 8462.691999999999
$ws.Range("K88").Value = 963.3333  # was 1144.2
$ws.Range("L88").Value = 9085.5  # was 8462.691999999999
$ws.Range("M88").Value = -557.3333  # was -738.2
$ws.Range("N88").Value = -9897.5  # was -9274.691999999999

$ws.Range("H91").Value = 6378.1113  # was 6429.778
$ws.Range("I91").Value = 963.3333  # was 1144.2
$ws.Range("J91").Value = 9085.5  # was 8462.691999999999
$ws.Range("K91").Value = 963.3333  # was 1144.2
$ws.Range("L91").Value = 9085.5  # was 8462.691999999999
$ws.Range("M91").Value = 440.6667  # was 259.8
$ws.Range("N91").Value = -11893.5  # was -11270.692

$ws.Range("H97").Value = 948.2727  # was 978.1
$ws.Range("I97").Value = 981.44446  # was 1022.875
$ws.Range("K97").Value = 981.44446  # was 1022.875
$ws.Range("M97").Value = -485.44446  # was -526.875

$ws.Range("H136").Value = 4846.7036  # was 4597.8623
$ws.Range("I136").Value = 5088.533  # was 4635.5884
$ws.Range("K136").Value = 15265.599  # was 13906.7652
$ws.Range("M136").Value = -12715.599  # was -11356.7652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 820  # was 708.4
$ws.Range("I4").Value = 244.5  # was 260.5
$ws.Range("J4").Value = 1395.5  # was 2500
$ws.Range("K4").Value = 244.5  # was 260.5
$ws.Range("L4").Value = 1395.5  # was 2500
$ws.Range("M4").Value = -129.5  # was -145.5
$ws.Range("N4").Value = -1625.5  # was -2730

$ws.Range("H20").Value = 2659  # was 2567.2222
$ws.Range("I20").Value = 2896.3  # was 2724.5454
$ws.Range("K20").Value = 2896.3  # was 2724.5454
$ws.Range("M20").Value = -2649.3  # was -2477.5454

$ws.Range("H54").Value = 19538  # was 19543.5
$ws.Range("J54").Value = 0  # was 19549
$ws.Range("L54").Value = 0  # was 19549
$ws.Range("N54").ClearContents()  # was -20517

$ws.Range("H64").Value = 5160.8  # was 4968.3
$ws.Range("I64").Value = 1291  # was 1217.6666
$ws.Range("J64").Value = 6128.25  # was 6575.7144
$ws.Range("K64").Value = 1291  # was 1217.6666
$ws.Range("L64").Value = 6128.25  # was 6575.7144
$ws.Range("M64").Value = -1066  # was -992.6666
$ws.Range("N64").Value = -6578.25  # was -7025.7144

$ws.Range("H67").Value = 5160.8  # was 4968.3
$ws.Range("I67").Value = 1291  # was 1217.6666
$ws.Range("J67").Value = 6128.25  # was 6575.7144
$ws.Range("K67").Value = 1291  # was 1217.6666
$ws.Range("L67").Value = 6128.25  # was 6575.7144
$ws.Range("M67").Value = -511  # was -437.6666
$ws.Range("N67").Value = -7688.25  # was -8135.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5013.3477  # was 5195.273
$ws.Range("I31").Value = 1471.5834  # was 1477
$ws.Range("J31").Value = 8877.091  # was 8913.546
$ws.Range("K31").Value = 1471.5834  # was 1477
$ws.Range("L31").Value = 8877.091  # was 8913.546
$ws.Range("M31").Value = -1176.5834  # was -1182
$ws.Range("N31").Value = -9467.091  # was -9503.546

$ws.Range("H34").Value = 5013.3477  # was 5195.273
$ws.Range("I34").Value = 1471.5834  # was 1477
$ws.Range("J34").Value = 8877.091  # was 8913.546
$ws.Range("K34").Value = 1471.5834  # was 1477
$ws.Range("L34").Value = 8877.091  # was 8913.546
$ws.Range("M34").Value = -1269.5834  # was -1275
$ws.Range("N34").Value = -9281.091  # was -9317.546

$ws.Range("H99").Value = 2007.1818  # was 2042.7778
$ws.Range("J99").Value = 2137.8  # was 2331.6667
$ws.Range("L99").Value = 2137.8  # was 2331.6667
$ws.Range("N99").Value = -5133.8  # was -5327.6667

$ws.Range("H126").Value = 2007.1818  # was 2042.7778
$ws.Range("J126").Value = 2137.8  # was 2331.6667
$ws.Range("L126").Value = 6413.400000000001  # was 6995.000100000001
$ws.Range("N126").Value = -11353.4  # was -11935.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8043.4  # was 7803.524
$ws.Range("J70").Value = 7790.273  # was 7391.5835
$ws.Range("L70").Value = 7790.273  # was 7391.5835
$ws.Range("N70").Value = -8330.273000000001  # was -7931.5835

$ws.Range("H73").Value = 8043.4  # was 7803.524
$ws.Range("J73").Value = 7790.273  # was 7391.5835
$ws.Range("L73").Value = 7790.273  # was 7391.5835
$ws.Range("N73").Value = -9662.273000000001  # was -9263.583500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1589.4166  # was 1593.9166
$ws.Range("J46").Value = 2200.8333  # was 2209.8333
$ws.Range("L46").Value = 2200.8333  # was 2209.8333
$ws.Range("N46").Value = -2576.8333  # was -2585.8333

$ws.Range("H68").Value = 2157.5293  # was 2641
$ws.Range("I68").Value = 2268.7  # was 2487.25
$ws.Range("J68").Value = 1998.7142  # was 2948.5
$ws.Range("K68").Value = 2268.7  # was 2487.25
$ws.Range("L68").Value = 1998.7142  # was 2948.5
$ws.Range("M68").Value = -1519.7  # was -1738.25
$ws.Range("N68").Value = -3496.7142  # was -4446.5

$ws.Range("H71").Value = 2157.5293  # was 2641
$ws.Range("I71").Value = 2268.7  # was 2487.25
$ws.Range("J71").Value = 1998.7142  # was 2948.5
$ws.Range("K71").Value = 11343.5  # was 12436.25
$ws.Range("L71").Value = 9993.571  # was 14742.5
$ws.Range("M71").Value = -7599.5  # was -8692.25
$ws.Range("N71").Value = -17481.571  # was -22230.5

$ws.Range("H92").Value = 0  # was 58500
$ws.Range("J92").Value = 0  # was 58500
$ws.Range("L92").Value = 0  # was 58500
$ws.Range("N92").ClearContents()  # was -63492
